$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "Senia Lucrezia"
$ws.Range("B29").Value = "Elia Battisti | U.SGUARNA"
$ws.Range("C29").Value = "Daniel Pedrotti | IMONTAGNA"
$ws.Range("D29").Value = "Michele Merighi | Clitoriders"
$ws.Range("E29").Value = "Gabriele Verona | CGB Gamberoni"
$ws.Range("F29").Value = "Gianni Sala | FC SALAGIARDINI"
